# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets
#
# 1. Insert a new "Player Info" sheet at the front of the workbook with the
#    player's basic info (ID, NAME, BATTING_HAND, BOWL_STYLE).
# 2. Rename the "MATCH_CARD_LINK" column to "MATCH_CODE" on both the
#    "ODI Batting" and "ODI Bowling" sheets, and replace the full scorecard
#    URL values with just the numeric match code.

$wb = $excel.ActiveWorkbook

$battingSheet = $wb.Worksheets.Item("ODI Batting")

# ---------------------------------------------------------------------
# 1. New "Player Info" sheet, inserted before "ODI Batting"
# ---------------------------------------------------------------------
$infoSheet = $wb.Worksheets.Add($battingSheet)
$infoSheet.Name = "Player Info"

# NOTE: sheet handles returned by Item()/Add() are index-anchored, not
# identity-anchored — inserting a sheet shifts every index after it, so
# re-resolve the handles we still need *after* the insert.
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# Copy the header formatting (bold, centered, bordered) from an existing
# header cell so the new sheet's header row matches the look of the others.
$battingSheet.Range("A1").Copy()
$infoSheet.Range("A1:D1").PasteSpecial(-4122)
[void]$infoSheet.Range("A1").Select()
$excel.CutCopyMode = $false

$infoSheet.Range("A1").Value = "ID"
$infoSheet.Range("B1").Value = "NAME"
$infoSheet.Range("C1").Value = "BATTING_HAND"
$infoSheet.Range("D1").Value = "BOWL_STYLE"

$infoSheet.Range("A2:D2").NumberFormat = "@"
$infoSheet.Range("A2").Value = "5954"
$infoSheet.Range("B2").Value = "Anderson Phillip"
$infoSheet.Range("C2").Value = "Right Handed"
$infoSheet.Range("D2").Value = "Right Arm Fast Medium"

# ---------------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK (column D) -> MATCH_CODE
# ---------------------------------------------------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingLinks = @{
    2 = "4451"
    3 = "4577"
    4 = "4580"
    5 = "4590"
    6 = "4606"
}
foreach ($row in $battingLinks.Keys) {
    $cell = $battingSheet.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $battingLinks[$row]
}

# ---------------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK (column B) -> MATCH_CODE
# ---------------------------------------------------------------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingLinks = @{
    2 = "4451"
    3 = "4577"
    4 = "4580"
    5 = "4590"
    6 = "4606"
}
foreach ($row in $bowlingLinks.Keys) {
    $cell = $bowlingSheet.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $bowlingLinks[$row]
}

[void]$infoSheet.Range("A1").Select()
